$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Item("Sheet3")

# Grab a reference to the "hyperlink" cell style (cellXf index 3 in the
# original file) from a cell that already carries it, so we can re-apply it
# later without Excel minting a brand new style record.
$hyperlinkStyle = $ws3.Range("A5").Style

# --- 1. Move the last data row of Sheet2 (A3:J3) to the end of Sheet3 (A7:J7) ---
$src = $ws2.Range("A3:J3")
$dst = $ws3.Range("A7:J7")
$dst.Value = $src.Value2
for ($i = 1; $i -le 10; $i++) {
  $dst.Cells.Item(1, $i).Style = $src.Cells.Item(1, $i).Style
}
# B column holds a numeric-looking string ("123456"); re-assert text storage.
$ws3.Range("B7").Value = "'123456"

# --- 2. Sheet2: retarget the surviving row's e-mail text, then drop row 3 ---
$ws2.Range("A2").Hyperlinks.Delete()
$ws2.Range("A2").Value = "mrsi.loris@Lfbnrm.caarrr"
$ws2.Range("A3:J3").EntireRow.Delete()

# Re-create the hyperlink on Sheet2!A2 (address text is unchanged).
$ws2.Hyperlinks.Add($ws2.Range("A2"), "mailto:mrsi.loris@Lfrm.ca") | Out-Null
$ws2.Range("A2").Style = $hyperlinkStyle

# --- 3. Sheet3: give the moved row its updated e-mail text + hyperlink ---
$ws3.Range("A7").Value = "Donefc.sat.ms@quetbhaeffseper.net"
$ws3.Hyperlinks.Add($ws3.Range("A7"), "mailto:Donefc.sat.ms@quetaeffseper.net") | Out-Null
$ws3.Range("A7").Style = $hyperlinkStyle

# --- 4. Fix up the selection shown on Sheet3 (without changing which tab
#        is active - Sheet2 stays the active tab, matching the original) ---
$ws3.Activate()
$ws3.Range("A7:XFD7").Select()
$ws2.Activate()
